$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.290.39'
$ws.Range("E2").Value = '  +1.69%  '
$ws.Range("D3").Value = '1.839.93'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("D4").Value = '''0.9995'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''243.07'
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("D6").Value = '''0.6853'
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").Value = '''1.000'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '''0.3026'
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").Value = '''0.07498'
$ws.Range("E9").Value = '  -1.73%  '
$ws.Range("D10").Value = '''23.19'
$ws.Range("E10").Value = '  +0.36%  '
$ws.Range("D11").Value = '''0.07648'
$ws.Range("E11").Value = '  -1.70%  '
$ws.Range("D12").Value = '1.843.88'
$ws.Range("E12").Value = '  +0.79%  '
$ws.Range("D13").Value = '''5.072'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = '''0.6844'
$ws.Range("E14").Value = '  +1.30%  '
$ws.Range("D15").Value = '''89.17'
$ws.Range("E15").Value = '  -3.80%  '
$ws.Range("D16").Value = '''6.285'
$ws.Range("E16").Value = '  -2.29%  '
$ws.Range("D17").Value = '29.258.83'
$ws.Range("E17").Value = '  +1.48%  '
$ws.Range("D18").Value = '''0.000008235'
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").Value = '2.090.00'
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("D20").Value = '''234.10'
$ws.Range("E20").Value = '  -2.92%  '
$ws.Range("D21").Value = '''12.57'
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '''7.454'
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("D24").Value = '''1.000'
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = '''0.1455'
$ws.Range("E25").Value = '  -1.57%  '
$ws.Range("E26").Value = '  -1.48%  '
$ws.Range("D27").Value = '''8.820'
$ws.Range("E27").Value = '  +1.32%  '
$ws.Range("D28").Value = '''18.07'
$ws.Range("E28").Value = '  -0.59%  '
$ws.Range("D29").Value = '''1.519'
$ws.Range("E29").Value = '  -1.28%  '
$ws.Range("D30").Value = '''4.215'
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("D31").Value = '''4.125'
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("D32").Value = '''1.201'
$ws.Range("E32").Value = '  +1.32%  '
$ws.Range("D33").Value = '''0.05126'
$ws.Range("E33").Value = '  +0.78%  '
$ws.Range("D34").Value = '''0.7684'
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("D36").Value = '''1.136'
$ws.Range("D37").Value = '''2.670'
$ws.Range("E37").Value = '  -0.83%  '
$ws.Range("D38").Value = '1.287.31'
$ws.Range("E38").Value = '  +3.81%  '
$ws.Range("D39").Value = '''0.01843'
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("E41").Value = '  -1.37%  '
$ws.Range("D42").Value = '''105.34'
$ws.Range("E42").Value = '  -1.33%  '
$ws.Range("D43").Value = '''1.000'
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '''5.634'
$ws.Range("E44").Value = '  -5.44%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''9.645'
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.990.85'
$ws.Range("E46").Value = '  +0.82%  '
$ws.Range("D47").Value = '''0.5204'
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("E48").Value = '  -1.05%  '
$ws.Range("D49").Value = '''1.758'
$ws.Range("E49").Value = '  +1.30%  '
$ws.Range("D50").Value = '''62.85'
$ws.Range("E50").Value = '  -1.04%  '
$ws.Range("D51").Value = '''0.05922'
$ws.Range("E51").Value = '  +0.88%  '
